$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '96.995.29'
$ws.Range("E2").Value = '  +0.28%  '

# Row 3
$ws.Range("D3").Value = '3.702.73'
$ws.Range("E3").Value = '  +0.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.56'
$ws.Range("E5").Value = '  -2.54%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.94'
$ws.Range("E6").Value = '  +2.59%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '656.92'
$ws.Range("E7").Value = '  -1.49%  '

# Row 8
$ws.Range("E8").Value = '  -0.13%  '

# Row 9
$ws.Range("E9").Value = '  -3.57%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("E10").Value = '  -0.03%  '

# Row 11
$ws.Range("D11").Value = '3.702.37'
$ws.Range("E11").Value = '  +0.55%  '

# Row 12
$ws.Range("E12").Value = '  -3.04%  '

# Row 13
$ws.Range("E13").Value = '  +1.49%  '

# Row 14
$ws.Range("E14").Value = '  +10.95%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.79'
$ws.Range("E15").Value = '  +2.78%  '

# Row 16
$ws.Range("D16").Value = '4.391.42'
$ws.Range("E16").Value = '  +0.50%  '

# Row 17
$ws.Range("D17").Value = '96.726.56'
$ws.Range("E17").Value = '  +0.26%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.95'
$ws.Range("E18").Value = '  -1.64%  '

# Row 19
$ws.Range("D19").Value = '3.697.94'
$ws.Range("E19").Value = '  +0.04%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.02'
$ws.Range("E20").Value = '  +0.67%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.70'
$ws.Range("E21").Value = '  +1.38%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.512'
$ws.Range("E22").Value = '  -4.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '523.18'
$ws.Range("E23").Value = '  +1.42%  '

# Row 24
$ws.Range("E24").Value = '  -1.39%  '

# Row 25
$ws.Range("E25").Value = '  +1.34%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.93'
$ws.Range("E26").Value = '  -0.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.04'
$ws.Range("E27").Value = '  -0.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.198'
$ws.Range("E28").Value = '  +17.95%  '

# Row 29
$ws.Range("E29").Value = '  +3.35%  '

# Row 30
$ws.Range("E30").Value = '  +1.65%  '

# Row 31
$ws.Range("E31").Value = '  -1.60%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.10%  '

# Row 33
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.189'
$ws.Range("E33").Value = '  +1.08%  '

# Row 34
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  +5.74%  '

# Row 35
$ws.Range("E35").Value = '  +0.42%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.33'
$ws.Range("E36").Value = '  -1.92%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '647.92'
$ws.Range("E37").Value = '  +5.43%  '

# Row 38
$ws.Range("E38").Value = '  +1.91%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.84'
$ws.Range("E39").Value = '  +1.28%  '

# Row 40
$ws.Range("E40").Value = '  +0.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.90'
$ws.Range("E41").Value = '  +11.15%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.05'
$ws.Range("E42").Value = '  +5.00%  '

# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.49'
$ws.Range("E43").Value = '  -4.90%  '

# Row 44
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.160'
$ws.Range("E44").Value = '  -0.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.960'
$ws.Range("E45").Value = '  +0.19%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0459'
$ws.Range("E46").Value = '  -0.46%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.447'
$ws.Range("E47").Value = '  +3.44%  '

# Row 48
$ws.Range("E48").Value = '  -0.77%  '

# Row 49
$ws.Range("E49").Value = '  +0.06%  '

# Row 50
$ws.Range("E50").Value = '  -0.89%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.53'
$ws.Range("E51").Value = '  +0.74%  '
